$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 118; existing rows 118-121 shift down to 119-122.
$ws.Rows.Item(118).Insert()

# Copy the style of the date cell (D) from the row above into new D118 so formatting matches.
$ws.Range("D117").Copy()
$ws.Range("D118").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new row 118 with the new weekly record.
$ws.Cells.Item(118, 1).Value = 8
$ws.Cells.Item(118, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(118, 3).Value = "Coquimbo"
$ws.Cells.Item(118, 4).Value = 44628
$ws.Cells.Item(118, 5).Value = 4
$ws.Cells.Item(118, 6).Value = 100112001
$ws.Cells.Item(118, 7).Value = "Berenjena"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 520
$ws.Cells.Item(118, 11).Value = 8500
$ws.Cells.Item(118, 12).Value = 9000
$ws.Cells.Item(118, 13).Value = 8750
$ws.Cells.Item(118, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(118, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(118, 16).Value = 175
$ws.Cells.Item(118, 17).Value = 50
$ws.Cells.Item(118, 18).Value = "Hortaliza"
